$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.45
$ws.Range("D2").Value = 1.28

$ws.Range("B3").Value = 1.38
$ws.Range("E3").Value = 1.3
$ws.Range("G3").Value = 0.71

$ws.Range("B4").Value = 1.43
$ws.Range("E4").Value = 1.18
$ws.Range("F4").Value = 0.96

$ws.Range("C5").Value = 1.36
$ws.Range("D5").Value = 1.39
$ws.Range("F5").Value = 1.08
$ws.Range("G5").Value = 0.54

$ws.Range("C7").Value = 2.14
$ws.Range("G7").Value = 1.2
